$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.1298123333333333
$ws.Range("H2").Value = 0.389437
$ws.Range("I2").Value = 0.01442185502613333
$ws.Range("J2").Value = 0.01442185502613333
$ws.Range("M2").Value = 200.005264
$ws.Range("N2").Value = 600.0157919999999
$ws.Range("O2").Value = 0.9465949791503665
$ws.Range("P2").Value = 0.9465949791503667
$ws.Range("Q2").Value = 25.96314999878933
$ws.Range("R2").Value = 233.668349989104
$ws.Range("S2").Value = 0.01365165555777229
$ws.Range("T2").Value = 0.0136516555577723
$ws.Range("G3").Value = 0.1298123333333333
$ws.Range("H3").Value = 0.389437
$ws.Range("I3").Value = 0.01442185502613333
$ws.Range("J3").Value = 0.01442185502613333
$ws.Range("O3").Value = 0.006425713585924051
$ws.Range("P3").Value = 0.006425713585924052
$ws.Range("Q3").Value = 0.1762440846985555
$ws.Range("R3").Value = 1.586196762287
$ws.Range("S3").Value = 0.00009267070977565201
$ws.Range("T3").Value = 0.00009267070977565204
$ws.Range("G4").Value = 0.1298123333333333
$ws.Range("H4").Value = 0.389437
$ws.Range("I4").Value = 0.01442185502613333
$ws.Range("J4").Value = 0.01442185502613333
$ws.Range("M4").Value = 9.926218666666667
$ws.Range("O4").Value = 0.04697930726370939
$ws.Range("P4").Value = 0.0469793072637094
$ws.Range("Q4").Value = 1.288545606296889
$ws.Range("R4").Value = 11.596910456672
$ws.Range("S4").Value = 0.0006775287585853893
$ws.Range("T4").Value = 0.0006775287585853896
$ws.Range("I5").Value = 0.7859600471098795
$ws.Range("J5").Value = 0.7859600471098797
$ws.Range("M5").Value = 200.005264
$ws.Range("N5").Value = 600.0157919999999
$ws.Range("O5").Value = 0.9465949791503665
$ws.Range("P5").Value = 0.9465949791503667
$ws.Range("Q5").Value = 1414.935773462731
$ws.Range("R5").Value = 12734.42196116457
$ws.Range("S5").Value = 0.7439858344069975
$ws.Range("T5").Value = 0.7439858344069977
$ws.Range("I6").Value = 0.7859600471098795
$ws.Range("J6").Value = 0.7859600471098797
$ws.Range("O6").Value = 0.006425713585924051
$ws.Range("P6").Value = 0.006425713585924052
$ws.Range("S6").Value = 0.00505035415270746
$ws.Range("T6").Value = 0.005050354152707461
$ws.Range("I7").Value = 0.7859600471098795
$ws.Range("J7").Value = 0.7859600471098797
$ws.Range("M7").Value = 9.926218666666667
$ws.Range("O7").Value = 0.04697930726370939
$ws.Range("P7").Value = 0.0469793072637094
$ws.Range("S7").Value = 0.03692385855017453
$ws.Range("T7").Value = 0.03692385855017455
$ws.Range("I8").Value = 0.1996180978639869
$ws.Range("J8").Value = 0.199618097863987
$ws.Range("M8").Value = 200.005264
$ws.Range("N8").Value = 600.0157919999999
$ws.Range("O8").Value = 0.9465949791503665
$ws.Range("P8").Value = 0.9465949791503667
$ws.Range("Q8").Value = 359.3653249130773
$ws.Range("R8").Value = 3234.287924217695
$ws.Range("S8").Value = 0.1889574891855965
$ws.Range("T8").Value = 0.1889574891855966
$ws.Range("I9").Value = 0.1996180978639869
$ws.Range("J9").Value = 0.199618097863987
$ws.Range("O9").Value = 0.006425713585924051
$ws.Range("P9").Value = 0.006425713585924052
$ws.Range("S9").Value = 0.001282688723440938
$ws.Range("T9").Value = 0.001282688723440938
$ws.Range("I10").Value = 0.1996180978639869
$ws.Range("J10").Value = 0.199618097863987
$ws.Range("M10").Value = 9.926218666666667
$ws.Range("O10").Value = 0.04697930726370939
$ws.Range("P10").Value = 0.0469793072637094
$ws.Range("S10").Value = 0.009377919954949453
$ws.Range("T10").Value = 0.009377919954949457